$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output")
$ws.Cells.Item(6, 2).Value = "latitude"
$ws.Cells.Item(7, 2).Value = "longitude"
